# "Se termino el desarrollo del servicio de carga de horarios"
# Rename the sheet, move the active selection, and swap the number formats
# that were applied to columns A (Código Carrera) and B (Nombre Comisión)
# on the "Importar" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Hoja1" -> "Importar"
$ws.Name = "Importar"

# Capture current (pre-swap) number formats used by column A and column B
$fmtA = $ws.Range("A2").NumberFormat
$fmtB = $ws.Range("B2").NumberFormat

# Swap number formats between column A (was currency) and column B (was time)
$ws.Range("A2:A4").NumberFormat = $fmtB
$ws.Range("B2:B4").NumberFormat = $fmtA

# Move the active selection to C7
$ws.Range("C7").Select()
